$wb = $excel.ActiveWorkbook

$wsPOM = $wb.Worksheets.Item("PredicateObjectMaps")

$wsPOM.Range("C8").Value = "iri"
$wsPOM.Range("C9").Value = "iri"
$wsPOM.Range("B10").Value = "recurso-trafico:tramot-{id_incidencia_nuevo}"
$wsPOM.Range("C10").Value = "iri"

$wsPOM.Range("B11").Select()
